$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "1.004", "0.4605").
# Force the cell to remain plain text (matching the original inlineStr cells)
# by temporarily applying a text number format, then resetting the style so the
# cell does not end up pinned to a non-default style index.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.718.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3857"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07844"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9810"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.891.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.991"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.695"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06955"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.717.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.265"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.097"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.081.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.878"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.986"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09324"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9164"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.293"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.334"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.322"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05777"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02074"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.644"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5619"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.746"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5288"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.833"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("E50").Value = "  +3.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "
